$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value = '35.509.86'
$ws.Cells.Item(2,5).Value = '  +2.97%  '

# Row 3
$ws.Cells.Item(3,4).Value = '1.840.34'
$ws.Cells.Item(3,5).Value = '  +1.93%  '

# Row 4
$ws.Cells.Item(4,5).Value = '  +0.31%  '

# Row 5
$ws.Cells.Item(5,4).NumberFormat = '@'
$ws.Cells.Item(5,4).Value = '231.65'
$ws.Cells.Item(5,4).Style = 'Normal'
$ws.Cells.Item(5,5).Value = '  +2.99%  '

# Row 6
$ws.Cells.Item(6,4).NumberFormat = '@'
$ws.Cells.Item(6,4).Value = '0.610'
$ws.Cells.Item(6,4).Style = 'Normal'
$ws.Cells.Item(6,5).Value = '  +1.09%  '

# Row 7
$ws.Cells.Item(7,5).Value = '  +0.30%  '

# Row 8
$ws.Cells.Item(8,4).NumberFormat = '@'
$ws.Cells.Item(8,4).Value = '43.85'
$ws.Cells.Item(8,4).Style = 'Normal'
$ws.Cells.Item(8,5).Value = '  +11.42%  '

# Row 9
$ws.Cells.Item(9,5).Value = '  +7.97%  '

# Row 10
$ws.Cells.Item(10,5).Value = '  +4.87%  '

# Row 11
$ws.Cells.Item(11,5).Value = '  +2.29%  '

# Row 12
$ws.Cells.Item(12,4).Value = '2.106.23'
$ws.Cells.Item(12,5).Value = '  +1.91%  '

# Row 13
$ws.Cells.Item(13,2).Value = 'Polygon'
$ws.Cells.Item(13,3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(13,4).NumberFormat = '@'
$ws.Cells.Item(13,4).Value = '0.676'
$ws.Cells.Item(13,4).Style = 'Normal'
$ws.Cells.Item(13,5).Value = '  +6.99%  '

# Row 14
$ws.Cells.Item(14,2).Value = 'WrappedEther'
$ws.Cells.Item(14,3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(14,4).Value = '1.830.47'
$ws.Cells.Item(14,5).Value = '  +1.48%  '

# Row 15
$ws.Cells.Item(15,5).Value = '  +2.07%  '

# Row 16
$ws.Cells.Item(16,5).Value = '  +8.21%  '

# Row 17
$ws.Cells.Item(17,4).Value = '35.475.93'
$ws.Cells.Item(17,5).Value = '  +2.94%  '

# Row 18
$ws.Cells.Item(18,4).NumberFormat = '@'
$ws.Cells.Item(18,4).Value = '70.37'
$ws.Cells.Item(18,4).Style = 'Normal'
$ws.Cells.Item(18,5).Value = '  +3.13%  '

# Row 19
$ws.Cells.Item(19,5).Value = '  +4.24%  '

# Row 20
$ws.Cells.Item(20,4).NumberFormat = '@'
$ws.Cells.Item(20,4).Value = '244.23'
$ws.Cells.Item(20,4).Style = 'Normal'
$ws.Cells.Item(20,5).Value = '  +1.71%  '

# Row 21
$ws.Cells.Item(21,4).NumberFormat = '@'
$ws.Cells.Item(21,4).Value = '12.03'
$ws.Cells.Item(21,4).Style = 'Normal'
$ws.Cells.Item(21,5).Value = '  +7.90%  '

# Row 22
$ws.Cells.Item(22,4).NumberFormat = '@'
$ws.Cells.Item(22,4).Value = '4.76'
$ws.Cells.Item(22,4).Style = 'Normal'
$ws.Cells.Item(22,5).Value = '  +16.18%  '

# Row 23
$ws.Cells.Item(23,5).Value = '  +0.30%  '

# Row 24
$ws.Cells.Item(24,4).NumberFormat = '@'
$ws.Cells.Item(24,4).Value = '2.23'
$ws.Cells.Item(24,4).Style = 'Normal'
$ws.Cells.Item(24,5).Value = '  +2.43%  '

# Row 25
$ws.Cells.Item(25,4).NumberFormat = '@'
$ws.Cells.Item(25,4).Value = '171.33'
$ws.Cells.Item(25,4).Style = 'Normal'
$ws.Cells.Item(25,5).Value = '  +0.11%  '

# Row 26
$ws.Cells.Item(26,4).NumberFormat = '@'
$ws.Cells.Item(26,4).Value = '7.96'
$ws.Cells.Item(26,4).Style = 'Normal'
$ws.Cells.Item(26,5).Value = '  +3.46%  '

# Row 27
$ws.Cells.Item(27,4).NumberFormat = '@'
$ws.Cells.Item(27,4).Value = '17.82'
$ws.Cells.Item(27,4).Style = 'Normal'
$ws.Cells.Item(27,5).Value = '  +0.87%  '

# Row 28
$ws.Cells.Item(28,4).NumberFormat = '@'
$ws.Cells.Item(28,4).Value = '0.122'
$ws.Cells.Item(28,4).Style = 'Normal'
$ws.Cells.Item(28,5).Value = '  -0.32%  '

# Row 29
$ws.Cells.Item(29,5).Value = '  +29.45%  '

# Row 30
$ws.Cells.Item(30,5).Value = '  +0.36%  '

# Row 31
$ws.Cells.Item(31,4).Value = '3.317.89'
$ws.Cells.Item(31,5).Value = '  +36.56%  '

# Row 32
$ws.Cells.Item(32,5).Value = '  +7.73%  '

# Row 33
$ws.Cells.Item(33,2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(33,3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(33,4).NumberFormat = '@'
$ws.Cells.Item(33,4).Value = '4.08'
$ws.Cells.Item(33,4).Style = 'Normal'
$ws.Cells.Item(33,5).Value = '  +6.20%  '

# Row 34
$ws.Cells.Item(34,2).Value = 'Filecoin'
$ws.Cells.Item(34,3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(34,4).NumberFormat = '@'
$ws.Cells.Item(34,4).Value = '3.94'
$ws.Cells.Item(34,4).Style = 'Normal'
$ws.Cells.Item(34,5).Value = '  +4.88%  '

# Row 35
$ws.Cells.Item(35,4).NumberFormat = '@'
$ws.Cells.Item(35,4).Value = '1.85'
$ws.Cells.Item(35,4).Style = 'Normal'
$ws.Cells.Item(35,5).Value = '  +1.88%  '

# Row 36
$ws.Cells.Item(36,4).NumberFormat = '@'
$ws.Cells.Item(36,4).Value = '95.55'
$ws.Cells.Item(36,4).Style = 'Normal'
$ws.Cells.Item(36,5).Value = '  +16.53%  '

# Row 37
$ws.Cells.Item(37,4).NumberFormat = '@'
$ws.Cells.Item(37,4).Value = '0.691'
$ws.Cells.Item(37,4).Style = 'Normal'
$ws.Cells.Item(37,5).Value = '  +7.72%  '

# Row 38
$ws.Cells.Item(38,4).NumberFormat = '@'
$ws.Cells.Item(38,4).Value = '1.13'
$ws.Cells.Item(38,4).Style = 'Normal'
$ws.Cells.Item(38,5).Value = '  +7.54%  '

# Row 39
$ws.Cells.Item(39,2).Value = 'Maker'
$ws.Cells.Item(39,3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(39,4).Value = '1.350.90'
$ws.Cells.Item(39,5).Value = '  +3.61%  '

# Row 40
$ws.Cells.Item(40,2).Value = 'InjectiveProtocol'
$ws.Cells.Item(40,3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(40,4).NumberFormat = '@'
$ws.Cells.Item(40,4).Value = '15.52'
$ws.Cells.Item(40,4).Style = 'Normal'
$ws.Cells.Item(40,5).Value = '  +11.86%  '

# Row 41
$ws.Cells.Item(41,5).Value = '  +5.98%  '

# Row 42
$ws.Cells.Item(42,5).Value = '  +4.80%  '

# Row 43
$ws.Cells.Item(43,5).Value = '  +6.50%  '

# Row 44
$ws.Cells.Item(44,4).NumberFormat = '@'
$ws.Cells.Item(44,4).Value = '1.28'
$ws.Cells.Item(44,4).Style = 'Normal'
$ws.Cells.Item(44,5).Value = '  +4.66%  '

# Row 45
$ws.Cells.Item(45,5).Value = '  +0.77%  '

# Row 46
$ws.Cells.Item(46,5).Value = '  +0.84%  '

# Row 47
$ws.Cells.Item(47,4).NumberFormat = '@'
$ws.Cells.Item(47,4).Value = '6.30'
$ws.Cells.Item(47,4).Style = 'Normal'
$ws.Cells.Item(47,5).Value = '  +9.62%  '

# Row 48
$ws.Cells.Item(48,5).Value = '  +0.92%  '

# Row 49
$ws.Cells.Item(49,4).Value = '2.015.61'
$ws.Cells.Item(49,5).Value = '  +2.46%  '

# Row 50
$ws.Cells.Item(50,5).Value = '  +0.36%  '

# Row 51
$ws.Cells.Item(51,4).NumberFormat = '@'
$ws.Cells.Item(51,4).Value = '103.36'
$ws.Cells.Item(51,4).Style = 'Normal'
$ws.Cells.Item(51,5).Value = '  +0.93%  '
